$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 8: Pagos (F8) and Inscrições homologadas (H8): 7 -> 8
$ws.Range("F8").Value = 8
$ws.Range("H8").Value = 8

# Row 15: Inscritos (E15): 151 -> 152
$ws.Range("E15").Value = 152

# Row 17: Inscritos (E17): 103 -> 104
$ws.Range("E17").Value = 104

# Row 36: Inscritos (E36): 92 -> 93
$ws.Range("E36").Value = 93

# Row 49: Inscritos (E49), Pagos (F49), Inscrições homologadas (H49)
$ws.Range("E49").Value = 62
$ws.Range("F49").Value = 35
$ws.Range("H49").Value = 35

# Row 72: Inscritos (E72): 37 -> 38
$ws.Range("E72").Value = 38

# Row 78: Inscritos (E78): 42 -> 43
$ws.Range("E78").Value = 43
